$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Update cell text values that changed (component version labels, etc.)
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "McPAT v0.7"
$ws.Range("A12").Value = "ORION v2.0 (May 2009)"
$ws.Range("A13").Value = "IntSim V1.0"
$ws.Range("A14").Value = "HotSpot v5.0"

# New row 17 - additional component "Iris"
$ws.Range("A17").Value = "Iris"

# Clear the stale "#code" software-location cell for [sst]disksim (C6)
$ws.Range("C6").ClearContents()

# ---------------------------------------------------------------------------
# 2. Rebuild the hyperlinks collection.
#    (The COM engine only supports wiping *all* hyperlinks at once via
#    Hyperlinks.Delete() - deleting a single Hyperlink object is a no-op -
#    so we clear everything and re-add every link we still need.)
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B9"), "http://www.cdkersey.com/qsim-web/")
$h = $ws.Hyperlinks.Add($ws.Range("C9"), "http://www.cdkersey.com/qsim-web/releases")
$h.TextToDisplay = "http://www.cdkersey.com/qsim-web/releases"
$ws.Range("C9").Value = "http://www.cdkersey.com/qsim-web/releases/"

$ws.Hyperlinks.Add($ws.Range("B2"), "http://sst.sandia.gov/")
$ws.Hyperlinks.Add($ws.Range("C2"), "http://code.google.com/p/sst-simulator/")
$ws.Hyperlinks.Add($ws.Range("C3"), "http://repo.gem5.org/")
$ws.Hyperlinks.Add($ws.Range("B3"), "http://www.m5sim.org/")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://wiki.umd.edu/DRAMSim2/")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://wiki.umd.edu/DRAMSim2/index.php?title=Main_Page", "Getting_DRAMSim2")
$ws.Hyperlinks.Add($ws.Range("B8"), "http://zesto.cc.gatech.edu/")
$ws.Hyperlinks.Add($ws.Range("C8"), "http://zesto.cc.gatech.edu/download.html")
$ws.Hyperlinks.Add($ws.Range("B6"), "http://www.pdl.cmu.edu/DiskSim/")
$ws.Hyperlinks.Add($ws.Range("C10"), "http://code.google.com/p/sst-simulator/")
$ws.Hyperlinks.Add($ws.Range("C11"), "http://code.google.com/p/sst-simulator/")
$ws.Hyperlinks.Add($ws.Range("C12"), "http://code.google.com/p/sst-simulator/")
$ws.Hyperlinks.Add($ws.Range("C13"), "http://code.google.com/p/sst-simulator/")
$ws.Hyperlinks.Add($ws.Range("C14"), "http://code.google.com/p/sst-simulator/")
$ws.Hyperlinks.Add($ws.Range("B11"), "http://www.hpl.hp.com/research/mcpat/")
$ws.Hyperlinks.Add($ws.Range("C5"), "http://code.google.com/p/sst-simulator/")

# New hyperlinks (to satisfy updated McPAT/ORION/IntSim/HotSpot entries).
# Hyperlinks.Add() does NOT rewrite the cell's displayed text by itself, so
# the new URL text has to be written into the cell explicitly as well.
$ws.Hyperlinks.Add($ws.Range("B12"), "http://projects.csail.mit.edu/cgi-bin/wiki/view/LSPgroup/OrionPage")
$ws.Range("B12").Value = "http://projects.csail.mit.edu/cgi-bin/wiki/view/LSPgroup/OrionPage"

$ws.Hyperlinks.Add($ws.Range("B13"), "http://deepaksekar.weebly.com/intsim.html")
$ws.Range("B13").Value = "http://deepaksekar.weebly.com/intsim.html"

$ws.Hyperlinks.Add($ws.Range("B14"), "http://lava.cs.virginia.edu/HotSpot/documentation.htm")
$ws.Range("B14").Value = "http://lava.cs.virginia.edu/HotSpot/documentation.htm"

# ---------------------------------------------------------------------------
# 3. Re-apply the "Hyperlink" cell style that Hyperlinks.Add() mangles, so
#    every cell ends up with the same visual style (font/wrap/valign) it
#    had (or should have) - this mirrors the original workbook's styling.
# ---------------------------------------------------------------------------
function Set-HyperlinkStyle($rangeAddr, $wrap, $centerV) {
    $r = $ws.Range($rangeAddr)
    $r.Style = "Hyperlink"
    # Only touch WrapText when it needs to be turned ON - explicitly setting
    # it to False forces the engine to mint a brand-new (unwanted) style.
    if ($wrap) { $r.WrapText = $true }
    if ($centerV) { $r.VerticalAlignment = -4108 }
}

Set-HyperlinkStyle "B9"  $true  $true
Set-HyperlinkStyle "C9"  $true  $true
Set-HyperlinkStyle "B2"  $true  $false
Set-HyperlinkStyle "C2"  $true  $false
Set-HyperlinkStyle "C3"  $true  $false
Set-HyperlinkStyle "B3"  $true  $false
Set-HyperlinkStyle "B4"  $false $false  # "Hyperlink" only - no wrap, no center
Set-HyperlinkStyle "C4"  $true  $false
Set-HyperlinkStyle "B8"  $true  $false
Set-HyperlinkStyle "C8"  $true  $false
Set-HyperlinkStyle "B6"  $true  $false
Set-HyperlinkStyle "C10" $true  $false
Set-HyperlinkStyle "C11" $true  $false
Set-HyperlinkStyle "C12" $true  $false
Set-HyperlinkStyle "C13" $true  $false
Set-HyperlinkStyle "C14" $true  $false
Set-HyperlinkStyle "B11" $true  $false
Set-HyperlinkStyle "C5"  $true  $false
Set-HyperlinkStyle "B12" $true  $false
Set-HyperlinkStyle "B13" $true  $false
Set-HyperlinkStyle "B14" $true  $false

# ---------------------------------------------------------------------------
# 4. Selection / active cell, matches the diff's sheetView/selection update.
# ---------------------------------------------------------------------------
$ws.Range("E20").Select()
